# "Generate Report for Handback" — refresh the localization-status report
# after a successful handback: the status text moves from "Ready for
# handoff" to "Handed back: in sync with en-US", the per-language handback
# timestamp advances, the (now resolved) handback "stale file" error detail
# is cleared, and the Status / Error Detail columns are re-sized so the
# new text fits.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn (col E) and de-de (col F) status columns
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Widen columns E & F so the longer status string fits (was ~17.22 chars).
$overview.Columns.Item(5).ColumnWidth = 29.166666666666664
$overview.Columns.Item(6).ColumnWidth = 29.166666666666664

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# Latest Handback DateTime refreshed to the new handback run.
$zhcn.Range("K2").Value = "2016-08-13 18:44:24"
$zhcn.Range("K3").Value = "2016-08-13 18:44:24"

# Error Detail cleared now that the handback file is in sync.
$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666664
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333332

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Latest Handback DateTime refreshed to the new handback run.
$dede.Range("K2").Value = "2016-08-13 18:44:33"
$dede.Range("K3").Value = "2016-08-13 18:44:33"

# Error Detail cleared now that the handback file is in sync.
$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.166666666666664
$dede.Columns.Item(16).ColumnWidth = 12.833333333333332
